$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the old "DoSetText" / "DoSendKeys" rows plus the following blank
# row (rows 20-22), shifting everything below up by three rows. The
# remaining two rows of the "Set Birthday" block (18-19, which used to be
# the duplicated "DoClick" rows) become the new "SetDate" function-call
# documentation rows.
$ws.Rows("20:22").Delete()

# Row 18: turn the old "Birthday / DoClick" row into the new
# "Functions / SetDate" action-doc row.
$ws.Range("C18").Value = "Functions"
$ws.Range("D18").Value = "SetDate"
$ws.Range("E18").Value = "field"
$ws.Range("F18").Value = "objectid"
$ws.Range("G18").Value = "Birthday"

# Row 19: turn it into the "Param" row describing the SetDate value.
$ws.Range("B19").Value = "Param"
$ws.Range("C19").ClearContents()
$ws.Range("D19").ClearContents()
$ws.Range("E19").Value = "value"
$ws.Range("F19").Value = "string"

# "12/25/1999" must stay a literal text value (matching the rest of the
# sheet, where this string already exists verbatim elsewhere) rather than
# be auto-converted into a date serial by Excel's input parsing. Enter it
# as a formula producing that text, then collapse the formula down to its
# computed value.
$ws.Range("G19").Formula = '="12/25/1999"'
$ws.Range("G19").Copy()
$ws.Range("G19").PasteSpecial(-4163)
$excel.CutCopyMode = 0
